# Gantry Parts List update
# -------------------------
# The "Control Board" line item (row 8), which used to link to the old
# Smoothieboard 4x product page, is replaced with the new part the lab
# actually orders from OpenBuilds: "Smoothieboard 5xC v1.1", at its new
# (higher) unit price. A stray leftover value that was sitting in E15
# (with nothing else around it) is also cleaned up.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 8: Smothieboard 4x item description / URL / price -> Smoothieboard 5xC v1.1
$ws.Range("B8").Value = "Smoothieboard 5xC v1.1"
$ws.Range("C8").Value = "http://openbuildspartstore.com/smoothieboard-5xc-v1-1/"
$ws.Range("E8").Value = 165.95

# Remove the stray leftover value in E15 (whole row is no longer used).
$ws.Range("E15").Clear() | Out-Null

# Match the saved selection / active cell.
$ws.Range("H14").Select()

$wb.Save()
